$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Diagnoses"

$headers = @("Id", "First Name", "Last Name", "Email", "Temperature", "Age", "Symptoms", "Total Ulhi", "Total Serious", "Total Common", "Total Less Common", "Current Fever", "Result")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
